# adj topic outline.
# The "Multiple Regression / Indicator Variables" topic (previously combined
# into a single week-11 row) is split across two weeks: "Multiple Regression"
# moves up into week 10's topic cell, and "Indicator Variables" remains alone
# in week 11's topic cell. The matching "Multiple Regression Assignment" eval
# item moves from week 11's eval cell up to week 10's eval cell, and the
# "Moderation and Confounding Assignment" heading is prefixed onto week 12's
# eval cell (which already listed the Poster/Peer-review Stage III items).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 10 row (row 11): "Flex time" -> "Multiple Regression"
$ws.Range("D11").Value = "Multiple Regression"

# Week 10 row (row 11): eval cell gets the Multiple Regression assignment note
$ws.Range("H11").Value = "Multiple Regression Assignment (Due 11/13)"

# Week 11 row (row 12): topic cell keeps only "Indicator Variables"
$ws.Range("D12").Value = "Indicator Variables"

# Week 11 row (row 12): eval cell no longer carries the assignment note
$ws.Range("H12").Value = ""

# Week 12 row (row 13): eval cell gains the Moderation and Confounding
# Assignment heading, prefixed before the existing Poster Prep / Peer Review
# Stage III lines.
$ws.Range("H13").Value = "Moderation and Confounding Assignment" + [char]10 + "* Poster Prep: Stage III (Due 11/15)" + [char]10 + "* Peer Review: Stage III (Due 11/17)"
